# Apply "Section_5 WIP, finished upto video# 45" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time_Series_Analysis_Forecastin")

# Fill in the new workout entry in row 19 (date, start time, end time)
$ws.Range("B19").Value = 44832
$ws.Range("C19").Value = 0.84375
$ws.Range("D19").Value = 0.92708333333333337

# E19 holds the elapsed-time formula, matching the pattern used by the rows above it
$ws.Range("E19").Formula = "=D19-C19"

# G19 gets the progress note for this session, using the same formatting as G18
$ws.Range("G18").Copy()
$ws.Range("G19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G19").Value = "Section 5, finished upto lesson# 45"

# Move the active selection to G20, as recorded in the saved view state
$ws.Range("G20").Select()

$wb.Save()
